$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (kNN)
$ws.Range("B2").Value = 0.02400603294372559
$ws.Range("C2").Value = 0.04416565895080567
$ws.Range("D2").Value = 0.01125965118408203
$ws.Range("E2").Value = 0.0294978141784668
$ws.Range("F2").Value = 0.003601789474487305
$ws.Range("G2").Value = 0.1209440231323242
$ws.Range("H2").Value = 0.02228279113769531
$ws.Range("I2").Value = 0.03478732109069824
$ws.Range("J2").Value = 0.01720924377441406
$ws.Range("K2").Value = 0.03222970962524414
$ws.Range("L2").Value = 0.005752182006835938
$ws.Range("M2").Value = 0.02220573425292969

# Row 3 (SVM)
$ws.Range("B3").Value = 0.09969134330749511
$ws.Range("C3").Value = 0.03379001617431641
$ws.Range("D3").Value = 0.01731934547424317
$ws.Range("E3").Value = 0.01303629875183105
$ws.Range("F3").Value = 0.008330821990966797
$ws.Range("G3").Value = 0.009400558471679688
$ws.Range("H3").Value = 0.1330410957336426
$ws.Range("I3").Value = 0.03856043815612793
$ws.Range("J3").Value = 0.08866133689880371
$ws.Range("K3").Value = 0.0279782772064209
$ws.Range("L3").Value = 0.02636990547180176
$ws.Range("M3").Value = 0.01359882354736328

# Row 4 (LR)
$ws.Range("B4").Value = 0.04292778968811035
$ws.Range("C4").Value = 0.02203035354614258
$ws.Range("D4").Value = 0.01440262794494629
$ws.Range("E4").Value = 0.01122441291809082
$ws.Range("F4").Value = 0.08300724029541015
$ws.Range("G4").Value = 0.01001482009887695
$ws.Range("H4").Value = 0.0386385440826416
$ws.Range("I4").Value = 0.02812857627868652
$ws.Range("J4").Value = 0.03425660133361817
$ws.Range("K4").Value = 0.02658829689025879
$ws.Range("L4").Value = 0.05747976303100586
$ws.Range("M4").Value = 0.01456212997436523

# Row 5 (NB)
$ws.Range("B5").Value = 0.03181700706481934
$ws.Range("C5").Value = 0.02595906257629394
$ws.Range("D5").Value = 0.02754673957824707
$ws.Range("E5").Value = 0.02444701194763184
$ws.Range("H5").Value = 0.03376898765563965
$ws.Range("I5").Value = 0.03006772994995117
$ws.Range("J5").Value = 0.02469053268432617
$ws.Range("K5").Value = 0.02999815940856933

# Row 6 (Ensemble)
$ws.Range("B6").Value = 0.4653768062591553
$ws.Range("C6").Value = 0.09450340270996094
$ws.Range("D6").Value = 0.543682336807251
$ws.Range("E6").Value = 0.1094675064086914
$ws.Range("F6").Value = 0.200335693359375
$ws.Range("G6").Value = 0.03403530120849609
$ws.Range("H6").Value = 1.090454530715942
$ws.Range("I6").Value = 0.21069016456604
$ws.Range("J6").Value = 0.7157362937927246
$ws.Range("K6").Value = 0.1745734214782715
$ws.Range("L6").Value = 0.3012603282928467
$ws.Range("M6").Value = 0.07853894233703614
